# Apply the change: append 3 new rows (168-170) to Sheet1 with the
# repeated Q/A pair about the maximum number of tracks in a single ODF file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$question = "What is the maximum number of tracks that can be specified within a single ODF file?"
$model = "llama3.2:latest"
$answer = "The maximum number of tracks that can be specified within a single ODF file is 200 tracks at most."

$startRow = 168
for ($i = 0; $i -lt 3; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $question
    $ws.Cells.Item($row, 2).Value = $model
    $ws.Cells.Item($row, 3).Value = $answer
}

$wb.Save()
